$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the values in D3:E18 while keeping their formatting/style
$ws.Range("D3:E18").ClearContents()

# Update the active selection to D3:E18 with D3 as the active cell
$ws.Range("D3:E18").Select()
